$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new values: P1=14, Q1=15
# (match the existing header formatting used by B1:O1 -- bold font,
#  thin box border, centered/top-aligned)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108
$ws.Range("P1:Q1").VerticalAlignment = -4160
$ws.Range("P1:Q1").Borders.LineStyle = 1

# For rows 2-25, swap values in columns I/K and M/O, and add P/Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P (new) = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q (new) = 2
}
